# "Generate Report for Handoff"
# The file "9490447e-e735-45e8-9da1-5aabfb41924c.md" has moved from
# "In Translation" to "Ready for handoff" for both locales (zh-cn, de-de).
# Update the per-locale worksheets (Status / Priority / Latest Handoff
# Datetime) and the roll-up "Overview" worksheet (locale status columns +
# Latest HO Xliff Generate Date), then autofit the columns whose displayed
# text just grew so the sheet matches how Excel would re-lay the report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: roll-up status + last handoff-xliff-generated date ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-18 06:13:48"
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("H3").Value = "2016-08-18 06:13:43"
$ws.Columns.Item(3).AutoFit()

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("H3").Value = "2016-08-18 06:13:48"
$ws.Columns.Item(3).AutoFit()
